# Aggiornato Avanzamento.xlsx da lbianco via Streamlit
#
# Data refresh: the "Produzione" column (D) for a batch of technicians in
# the September progress sheet was re-pulled from the source system, which
# changed 23 values. Column F ("Avanzamento") is a shared formula
# (=Dn-(Dn*En)/100) so it recomputes automatically once D is written.
# The sheet's scroll position / active selection is also updated to match
# where the author left the view (top-left A63, active cell G69).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D67").Value  = 43.242674199327503
$ws.Range("D68").Value  = 26.968554675444
$ws.Range("D72").Value  = 32.000048780487802
$ws.Range("D76").Value  = 53.331171415930001
$ws.Range("D79").Value  = 36.556056285178201
$ws.Range("D80").Value  = 35.585638888888901
$ws.Range("D83").Value  = 31.4012015873016
$ws.Range("D85").Value  = 31.2606573391813
$ws.Range("D87").Value  = 40.960740823549003
$ws.Range("D88").Value  = 36.235617341826
$ws.Range("D92").Value  = 25.3093183413974
$ws.Range("D96").Value  = 43.757928379987
$ws.Range("D98").Value  = 32.315026785714302
$ws.Range("D102").Value = 26.0163981880644
$ws.Range("D103").Value = 29.964407982952501
$ws.Range("D104").Value = 39.553742521367603
$ws.Range("D105").Value = 36.801722222222203
$ws.Range("D115").Value = 33.464523090277801
$ws.Range("D117").Value = 33.310427835648198
$ws.Range("D124").Value = 30.4467243076515
$ws.Range("D126").Value = 70.316847293609598
$ws.Range("D127").Value = 35.267092948718002
$ws.Range("D129").Value = 30.063688741721901

# Move the view/selection to match the author's final cursor position
# (topLeftCell A63 / activeCell G69 in the saved sheetView).
$ws.Range("A63").Select() | Out-Null
$ws.Range("G69").Select() | Out-Null
